$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("B6").Value = 3726
$ws.Range("C6").Value = 30.798479087452474
$ws.Range("D6").Value = 3726
$ws.Range("E6").Value = 30.798479087452474
$ws.Range("F6").Value = 3726
$ws.Range("G6").Value = 30.798479087452474

# Row 7
$ws.Range("B7").Value = 3963
$ws.Range("C7").Value = 32.757480575301699
$ws.Range("D7").Value = 3963
$ws.Range("E7").Value = 32.757480575301699
$ws.Range("F7").Value = 3963
$ws.Range("G7").Value = 32.757480575301699

# Row 8
$ws.Range("B8").Value = 8098
$ws.Range("C8").Value = 66.936683749380066
$ws.Range("D8").Value = 8098
$ws.Range("E8").Value = 66.936683749380066
$ws.Range("F8").Value = 8098
$ws.Range("G8").Value = 66.936683749380066

# Row 17 (only B,C,D,E change; F,G remain as-is)
$ws.Range("B17").Value = 33
$ws.Range("C17").Value = 0.56065239551478085
$ws.Range("D17").Value = 33
$ws.Range("E17").Value = 0.56065239551478085

# Row 18 (only B,C,D,E change; F,G remain as-is)
$ws.Range("B18").Value = 5796
$ws.Range("C18").Value = 98.470948012232412
$ws.Range("D18").Value = 5796
$ws.Range("E18").Value = 98.470948012232412
